$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22/23: Toncoin and Uniswap swap positions
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").Value = "2.54"
$ws.Range("E22").Value = "  +4.85%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "4.41"
$ws.Range("E23").Value = "  -0.15%  "

# Remaining D/E updates
$ws.Range("D2").Value = "27.165.81"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "1.636.95"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "216.76"
$ws.Range("E5").Value = "  -0.45%  "
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.253"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "20.02"
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").Value = "1.865.51"
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").Value = "1.635.70"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "66.34"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "27.160.56"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "216.45"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D24").Value = "9.12"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("D25").Value = "147.78"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("D29").Value = "15.64"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("E33").Value = "  +0.43%  "
$ws.Range("D34").Value = "1.299.68"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("D37").Value = "0.0176"
$ws.Range("E37").Value = "  -0.33%  "
$ws.Range("D38").Value = "0.549"
$ws.Range("E38").Value = "  +2.40%  "
$ws.Range("D39").Value = "0.854"
$ws.Range("E39").Value = "  +1.58%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "2.24"
$ws.Range("E42").Value = "  +5.64%  "
$ws.Range("D43").Value = "5.34"
$ws.Range("E43").Value = "  -0.60%  "
$ws.Range("D44").Value = "1.775.87"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "61.91"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "91.28"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").Value = "0.0514"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("E51").Value = "  -0.75%  "
